$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.855.86"
$ws.Range("E2").Value = "  +0.21%  "
$ws.Range("D3").Value = "1.636.66"
$ws.Range("E3").Value = "  -0.13%  "
$ws.Range("E4").Value = "  -0.44%  "
$ws.Range("D5").Formula = '="216.76"'
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = "  -0.89%  "
$ws.Range("D6").Formula = '="0.510"'
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Value = "  +1.69%  "
$ws.Range("E7").Value = "  -0.45%  "
$ws.Range("D8").Formula = '="0.256"'
$ws.Range("D8").Copy()
$ws.Range("D8").PasteSpecial(-4163)
$ws.Range("E8").Value = "  +1.97%  "
$ws.Range("E9").Value = "  +0.13%  "
$ws.Range("D10").Formula = '="19.91"'
$ws.Range("D10").Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("E10").Value = "  +3.22%  "
$ws.Range("D11").Formula = '="0.0845"'
$ws.Range("D11").Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("E11").Value = "  +0.11%  "
$ws.Range("D12").Value = "1.865.91"
$ws.Range("E12").Value = "  -0.11%  "
$ws.Range("D13").Value = "1.640.16"
$ws.Range("E13").Value = "  +0.05%  "
$ws.Range("D14").Formula = '="4.11"'
$ws.Range("D14").Copy()
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range("E14").Value = "  -0.74%  "
$ws.Range("D15").Formula = '="0.528"'
$ws.Range("D15").Copy()
$ws.Range("D15").PasteSpecial(-4163)
$ws.Range("E15").Value = "  +0.51%  "
$ws.Range("D16").Formula = '="66.69"'
$ws.Range("D16").Copy()
$ws.Range("D16").PasteSpecial(-4163)
$ws.Range("E16").Value = "  +2.83%  "
$ws.Range("D17").Value = "26.853.49"
$ws.Range("E17").Value = "  +0.24%  "
$ws.Range("E18").Value = "  -0.78%  "
$ws.Range("D19").Formula = '="219.28"'
$ws.Range("D19").Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("E19").Value = "  +1.48%  "
$ws.Range("E20").Value = "  -0.54%  "
$ws.Range("D21").Formula = '="6.74"'
$ws.Range("D21").Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("E21").Value = "  +3.17%  "
$ws.Range("D22").Formula = '="4.38"'
$ws.Range("D22").Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("E22").Value = "  +0.76%  "
$ws.Range("E23").Value = "  +3.18%  "
$ws.Range("D24").Formula = '="9.15"'
$ws.Range("D24").Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("E24").Value = "  +0.00%  "
$ws.Range("D25").Formula = '="147.02"'
$ws.Range("D25").Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("E25").Value = "  -0.15%  "
$ws.Range("E26").Value = "  -0.55%  "
$ws.Range("E27").Value = "  +4.43%  "
$ws.Range("E28").Value = "  +0.38%  "
$ws.Range("D29").Formula = '="15.78"'
$ws.Range("D29").Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("E29").Value = "  +0.42%  "
$ws.Range("E30").Value = "  -0.42%  "
$ws.Range("E31").Value = "  -1.31%  "
$ws.Range("E32").Value = "  -1.38%  "
$ws.Range("E33").Value = "  +0.83%  "
$ws.Range("E34").Value = "  +0.87%  "
$ws.Range("D35").Value = "1.255.97"
$ws.Range("E35").Value = "  -0.48%  "
$ws.Range("E36").Value = "  -0.24%  "
$ws.Range("E37").Value = "  +2.04%  "
$ws.Range("E38").Value = "  +0.84%  "
$ws.Range("E39").Value = "  +2.03%  "
$ws.Range("E40").Value = "  -0.45%  "
$ws.Range("E41").Value = "  +0.42%  "
$ws.Range("D42").Formula = '="5.41"'
$ws.Range("D42").Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("E42").Value = "  +1.41%  "
$ws.Range("D43").Value = "1.780.62"
$ws.Range("E43").Value = "  +0.14%  "
$ws.Range("E44").Value = "  -1.63%  "
$ws.Range("D45").Formula = '="61.55"'
$ws.Range("D45").Copy()
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range("E45").Value = "  +2.22%  "
$ws.Range("D46").Formula = '="91.66"'
$ws.Range("D46").Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("E46").Value = "  -0.54%  "
$ws.Range("E47").Value = "  +0.60%  "
$ws.Range("E48").Value = "  +0.53%  "
$ws.Range("E49").Value = "  -0.44%  "
$ws.Range("D50").Formula = '="7.65"'
$ws.Range("D50").Copy()
$ws.Range("D50").PasteSpecial(-4163)
$ws.Range("E50").Value = "  +1.77%  "
$ws.Range("D51").Formula = '="0.0959"'
$ws.Range("D51").Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range("E51").Value = "  -0.28%  "
